$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.861.41"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.627.18"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'214.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.500"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.0630"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.852.46"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "1.619.86"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'0.543"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "0.0₃0756"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'62.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "25.851.71"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'192.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'142.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'0.0496"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "1.129.80"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'99.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'0.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "1.763.30"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'56.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0530"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.415"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.05%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0958"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.65%  "
